# Mark the test-data row as "Used" in the Is_Used column (G) of the
# New_Registration sheet, in line with TestNG setup for parallel test
# execution (tracking which rows of test data have already been consumed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New_Registration")

$ws.Range("G2").Value = "Used"
